$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text change
$ws.Range("H1").Value = "Net migration"

# Row 4
$ws.Range("H4").Value = -29158.12
$ws.Range("I4").Value = 0.19
$ws.Range("J4").Value = -1.43

# Row 5
$ws.Range("H5").Value = 1099373.6
$ws.Range("I5").Value = -1.11
$ws.Range("J5").Value = 0.39

# Row 6 (style changes from s=3 to s=5; B9 already has the target style)
$ws.Range("B9").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("H6").Value = -698788.4399999999
$ws.Range("I6").Value = 0.45
$ws.Range("J6").Value = -1.19

# Row 7
$ws.Range("H7").Value = 409604.44
$ws.Range("I7").Value = -0.21
$ws.Range("J7").Value = 1.18

# Row 8 (style changes from s=5 to s=4; E5 already has the target style)
$ws.Range("E5").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("H8").Value = 1515907.4
$ws.Range("I8").Value = -0.09
$ws.Range("J8").Value = 0.05

# Row 9 (style changes from s=4 to s=3; H4 already has the target style)
$ws.Range("H4").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("H9").Value = -539396.36
$ws.Range("I9").Value = 0.35
$ws.Range("J9").Value = 1.73
